$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.537.70'
$ws.Range('D3').Value = '1.674.87'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '219.98'
$ws.Range('D6').Value = '0.5297'
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('D9').Value = '0.06393'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').Value = '21.84'
$ws.Range('E10').Value = '  +4.84%  '
$ws.Range('D11').Value = '0.07810'
$ws.Range('E11').Value = '  +1.90%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.675.91'
$ws.Range('E12').Value = '  +1.65%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.497'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').Value = '0.5588'
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').Value = '0.0₅8344'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '65.76'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').Value = '26.562.59'
$ws.Range('E17').Value = '  +1.90%  '
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = '4.781'
$ws.Range('E19').Value = '  +1.16%  '
$ws.Range('D20').Value = '193.41'
$ws.Range('E20').Value = '  +2.60%  '
$ws.Range('E21').Value = '  +1.56%  '
$ws.Range('D22').Value = '6.322'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '0.1273'
$ws.Range('E24').Value = '  +4.39%  '
$ws.Range('D25').Value = '138.72'
$ws.Range('E25').Value = '  -5.10%  '
$ws.Range('D26').Value = '7.407'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  +3.24%  '
$ws.Range('E28').Value = '  +3.06%  '
$ws.Range('D29').Value = '0.06277'
$ws.Range('E29').Value = '  +5.56%  '
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').Value = '3.606'
$ws.Range('E31').Value = '  +6.16%  '
$ws.Range('D32').Value = '3.424'
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('D33').Value = '1.694'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('D34').Value = '1.012'
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('D35').Value = '0.6199'
$ws.Range('E35').Value = '  +10.30%  '
$ws.Range('D36').Value = '2.421'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '2.785'
$ws.Range('E37').Value = '  +1.10%  '
$ws.Range('D38').Value = '0.01622'
$ws.Range('E38').Value = '  +0.70%  '
$ws.Range('D39').Value = '6.081'
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('D40').Value = '1.096.51'
$ws.Range('E40').Value = '  +6.67%  '
$ws.Range('D41').Value = '0.8603'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').Value = '100.62'
$ws.Range('E43').Value = '  +1.39%  '
$ws.Range('D44').Value = '1.821.13'
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').Value = '58.98'
$ws.Range('E45').Value = '  +5.70%  '
$ws.Range('D46').Value = '0.0₈108'
$ws.Range('E46').Value = '  -3.67%  '
$ws.Range('D47').Value = '8.198'
$ws.Range('E47').Value = '  +1.39%  '
$ws.Range('D48').Value = '1.535'
$ws.Range('E48').Value = '  +11.09%  '
$ws.Range('D49').Value = '0.9996'
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('D50').Value = '0.05193'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').Value = '6.023'
$ws.Range('E51').Value = '  +1.47%  '
